$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (dates in column A, values in column B)
$dates = @(
    "2024-11-15",
    "2023-07-09",
    "2023-08-13",
    "2023-06-27",
    "2023-03-14",
    "2023-03-09",
    "2023-03-07",
    "2023-03-06",
    "2023-03-05",
    "2023-03-04",
    "2023-02-27",
    "2023-02-24",
    "2023-02-23",
    "2023-02-16",
    "2023-02-05",
    "2022-12-13",
    "2022-12-11",
    "2022-12-07",
    "2024-11-16",
    "2024-11-17",
    "2024-11-18",
    "2024-11-20",
    "2024-11-19",
    "2024-11-21",
    "2024-11-22",
    "2024-11-23",
    "2024-11-24"
)

$values = @(
    0.1764,
    0.08015,
    0.07742,
    0.07378,
    0.06716,
    0.06533,
    0.0674,
    0.06723,
    0.06673999999999999,
    0.06714000000000001,
    0.06955,
    0.06956,
    0.06943000000000001,
    0.07087,
    0.06381000000000001,
    0.05353,
    0.05481,
    0.05365,
    0.1908,
    0.2,
    0.1992,
    0.1994,
    0.202,
    0.1947,
    0.1986,
    0.2046,
    0.2119
)

$startRow = 995
$endRow = $startRow + $dates.Count - 1

# Force column A for the new rows to be stored as text (not auto-parsed as dates
# by Excel's smart input detection), matching the existing "t=s" shared-string
# cells used throughout the rest of column A.
$dateRange = $ws.Range("A$startRow`:A$endRow")
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $dates.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Restore default (unstyled) cell formatting so the new cells match the
# plain, unstyled cells used by the rest of the data rows.
$dateRange.Style = "Normal"
